$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.39"
$ws.Range("E2").Value = "'0.30%"
$ws.Range("D3").Value = "'41.00"
$ws.Range("E3").Value = "'2.59%"
$ws.Range("D4").Value = "'5.125"
$ws.Range("E4").Value = "'0.19%"
$ws.Range("D5").Value = "'0.07626"
$ws.Range("E5").Value = "'-1.34%"
$ws.Range("D6").Value = "'1.606"
$ws.Range("E6").Value = "'-0.19%"
$ws.Range("E7").Value = "'1.90%"
$ws.Range("D8").Value = "'0.9041"
$ws.Range("E8").Value = "'2.19%"
$ws.Range("D9").Value = "'0.1120"
$ws.Range("E9").Value = "'12.66%"
$ws.Range("D10").Value = "'0.1778"
$ws.Range("E10").Value = "'2.13%"
$ws.Range("D11").Value = "'0.09140"
$ws.Range("E11").Value = "'1.19%"
$ws.Range("D12").Value = "'0.04278"
$ws.Range("E12").Value = "'-3.40%"
$ws.Range("E13").Value = "'-0.34%"
$ws.Range("D14").Value = "'0.001254"
$ws.Range("E14").Value = "'-0.94%"
$ws.Range("D15").Value = "'0.005724"
$ws.Range("E15").Value = "'-1.95%"
$ws.Range("E16").Value = "'-0.11%"
$ws.Range("D17").Value = "'4.250"
$ws.Range("E17").Value = "'0.22%"
$ws.Range("E18").Value = "'0.64%"
$ws.Range("D19").Value = "'6.663"
$ws.Range("E19").Value = "'-6.44%"
$ws.Range("E20").Value = "'1.21%"
$ws.Range("D21").Value = "'0.2802"
$ws.Range("E21").Value = "'-1.63%"
$ws.Range("D22").Value = "'0.04069"
$ws.Range("E22").Value = "'-1.01%"
$ws.Range("D23").Value = "'0.001242"
$ws.Range("E23").Value = "'3.79%"
$ws.Range("E24").Value = "'0.84%"
$ws.Range("D26").Value = "'0.0003748"
$ws.Range("D38").Value = "'0.02382"
$ws.Range("E38").Value = "'1.42%"
$ws.Range("D39").Value = "'0.05197"
$ws.Range("E39").Value = "'-0.05%"
$ws.Range("D40").Value = "'0.007778"
$ws.Range("E40").Value = "'-1.83%"
$ws.Range("D41").Value = "'0.1303"
$ws.Range("E41").Value = "'-1.54%"
$ws.Range("E42").Value = "'13.31%"
$ws.Range("D44").Value = "'0.007941"
$ws.Range("E44").Value = "'-9.36%"
$ws.Range("D45").Value = "'0.3085"
$ws.Range("E45").Value = "'-7.22%"
$ws.Range("D46").Value = "'0.00007008"
$ws.Range("E46").Value = "'6.82%"
$ws.Range("E47").Value = "'-0.07%"
$ws.Range("D48").Value = "'0.03177"
$ws.Range("E48").Value = "'786.59%"
$ws.Range("E50").Value = "'-0.07%"
$ws.Range("E51").Value = "'-0.07%"
